$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1
$ws.Range("H1").Value = "Save"

# Copy formatting from the neighboring "sum" header (G1) so the new
# "Save" header matches the existing header style (bold, bordered,
# centered) used by the other column headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the Save column values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
